$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(1, 1).Value2 = 'Input'
$ws.Cells.Item(1, 2).Value2 = 'Chatbot Response'

$ws.Cells.Item(2, 1).Value2 = 'What is a chatbot?'
$ws.Cells.Item(2, 2).Value2 = 'A chatbot is a software application used to conduct an on-line chat conversation via text or text-to-speech, in lieu of providing direct contact with a live human agent. 
'

$ws.Cells.Item(3, 1).Value2 = 'Does Mega offer a chatbot?'
$ws.Cells.Item(3, 2).Value2 = 'Mega is currently developing a chatbot as part of MegaACD. 
'

$ws.Cells.Item(4, 1).Value2 = 'What will Mega chatbot do?
Mega chatbot applications
Mega chatbot uses
Mega chatbot purpose
Mega chatbot detail
more about mega chatbot
details'
$ws.Cells.Item(4, 2).Value2 = 'Mega chatbot is an AI-based chatbot that aims to lower the dependency on human agents by handling many of the everyday interactions that take place between customers and agents in a contact center setting. 

With the help of Mega chatbot, our clients will be able to:

- Provide information about their products and services
- Help customers make a better choice when shopping for products
- Resolve problems that customers face in  everyday lives when using their products.
'

$ws.Cells.Item(5, 1).Value2 = 'What features are currently available in Mega chatbot?
Mega chatbot features'
$ws.Cells.Item(5, 2).Value2 = 'Mega chatbot currently has following features:

- Ability to respond correctly to questions similar to the ones provided in data, even when asked differently
- A web interface provided to easily upload data for chatbot''s training
'

$ws.Cells.Item(6, 1).Value2 = 'How can data be provided to the chatbot for training?
chatbot training data'
$ws.Cells.Item(6, 2).Value2 = 'A web interface has been provided where data (all possible questions and answers in an Excel file) can be uploaded easily for the chatbot to consume. 
'

$ws.Cells.Item(7, 1).Value2 = 'How will the chabot know how to respond to questions asked?
How chatbot knows answers
How chatbot respond'
$ws.Cells.Item(7, 2).Value2 = 'Mega chatbot uses neural network and probability algorithm to identify the intent of the user input and present the most accurate response based on the data provided.

When the chabot is provided data with any number of questions and their corresponding responses, it trains itself to learn all the questions and their responses. Once training is complete, when a user asks the chatbot a question, chatbot will identify the intent of the question asked by comparing it to the questions it was trained on using keywords and provide the most appropriate response.
'

$ws.Cells.Item(8, 1).Value2 = 'How will the chatbot be trained?
Chatbot training
chatbot file upload
chatbot excel
chatbot upload data'
$ws.Cells.Item(8, 2).Value2 = 'A web interface has been provided on which data (in Excel files) can be uploaded. Once the upload is complete, chatbot will pick up the file and begin the training process.
'

$ws.Cells.Item(9, 1).Value2 = 'What happens during the training process?
Training process
how chatbot is trained
how chatbot learn'
$ws.Cells.Item(9, 2).Value2 = 'When the chatbot trains, it reads through all the questions provided in the data, identifies keywords and looks up their synonyms. It then uses the frequency of the key words in the questions and their synonyms to associate them to responses provided in data. This enables chatbot to understand questions that asked are differently from the data questions and respond to them.
'

$ws.Cells.Item(10, 1).Value2 = 'What will happen if chatbot is asked a question that it does not understand but was trained on related question?
Unrelated question'
$ws.Cells.Item(10, 2).Value2 = 'It is possible that the chatbot is asked a question that it does not undertand even though it was trained on a question that was worded/phrased differently but had the same meaning. In this case it will respond by saying "Sorry, I do not understad the question. Can you please try saying it a different way?
'

$ws.Cells.Item(11, 1).Value2 = 'What will happen if chatbot is asked a question it was not trained on?
Response for question not trained'
$ws.Cells.Item(11, 2).Value2 = 'In this case the chatbot will respond by saying "Sorry, I do not understad the question. Can you please try saying it a different way?
'

$ws.Cells.Item(12, 1).Value2 = 'What are  the current limitations of the chatbot?
Limitation of chatbot'
$ws.Cells.Item(12, 2).Value2 = 'Mega chatbot''s current limitations are:

- It is not self-learning and data has to be provided manually via a web interface
- It can only respond to questions that have been provided in the data 
- It treats each question asked as a separate question and does not understand the context in which they are being asked. It is therefore not able to have a back and forth conversation on a specific topic
- It is not able to access to customer data and therefore cannot provide support specific to individual customers
'

$ws.Cells.Item(13, 1).Value2 = 'What features are planned to be available in the chatbot in the future?
Future features in our chatbot
coming features in chatbot'
$ws.Cells.Item(13, 2).Value2 = 'Mega chatbot will have the following features:

- Use neural network and deep learning to understand the context of customers'' queries in addition to the intent
- Ability to access customer data through CRM and billing systems after customers log in, to identify products used and customers'' history with the company
- Provide customer an option (through a graphical button or text input) to speak to a human agent at any point in the chat
- Ability to monitor real-time chat between human agents and customers and use the chat as data to train itself and enhance its capability to better interact with customers
- Calculate the probability of being able to respond accurately and appropritely to all user inputs and only respond if probability is high. A low probability calculation will lead to customer being connected to an agent. This feature will help maintain good customer experience when using chatbot
'

$ws.Cells.Item(14, 1).Value2 = 'What is the benefit or  purpose of using neural network and deep learning in Mega chatbot?
Neural networks in chatbot
Deep learning in chatbots
purpose of Neural networks
benefit of deep learning'
$ws.Cells.Item(14, 2).Value2 = 'Using neural network and deep learning will enable Mega chatbot to understand the intent of user''s input and the context of the conversation. This will enable the chatbot to:

1) Respond to questions correctly even when asked differently from the questions in data provided
2) Understand follow up  questions users may ask in reply to chatbot''s response'

$ws.Cells.Item(15, 1).Value2 = 'What is the benefit or  purpose of allowing Mega chatbot to access customer data?
why access customer data'
$ws.Cells.Item(15, 2).Value2 = 'Through access to customer data through CRM and billing systems, products used by customers and their history with a company can be identified.  This will help better understand the context of customers'' queries and provide support specific to individual customers based on their needs. 
'

$ws.Cells.Item(16, 1).Value2 = 'What is out of scope for the chatbot?
Chatbot out of scope'
$ws.Cells.Item(16, 2).Value2 = 'Mega chatbot will not be a "chatty" bot and will not respond to customer inputs that are not related to a company''s business, services, or products. 
'

$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("A20").Select()
